$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 2000.5
$ws.Range("I2").Value = 1750.0
$ws.Range("J2").Value = 2251.0
$ws.Range("K2").Value = 1750.0
$ws.Range("L2").Value = 2251.0
$ws.Range("M2").Value = -1637.0
$ws.Range("N2").Value = -2477.0
# Row 19
$ws.Range("H19").Value = 2180.2
$ws.Range("I19").Value = 1518.0
$ws.Range("K19").Value = 1518.0
$ws.Range("M19").Value = -1343.0
# Row 69
$ws.Range("H69").Value = 3006.0
$ws.Range("I69").Value = 3006.0
$ws.Range("K69").Value = 9018.0
$ws.Range("M69").Value = -8144.0
# Row 72
$ws.Range("H72").Value = 3006.0
$ws.Range("I72").Value = 3006.0
$ws.Range("K72").Value = 27054.0
$ws.Range("M72").Value = -22686.0
# Row 92
$ws.Range("H92").Value = 327.58334
$ws.Range("I92").Value = 296.9091
$ws.Range("K92").Value = 296.9091
$ws.Range("M92").Value = 951.0908999999999
# Row 98
$ws.Range("H98").Value = 731.875
$ws.Range("I98").Value = 309.16666
$ws.Range("K98").Value = 309.16666
$ws.Range("M98").Value = 1188.83334
# Row 99
$ws.Range("H99").Value = 326.0
$ws.Range("I99").Value = 326.0
$ws.Range("K99").Value = 978.0
$ws.Range("M99").Value = 520.0
# Row 115
$ws.Range("H115").Value = 5539.2856
$ws.Range("I115").Value = 5539.2856
$ws.Range("K115").Value = 16617.8568
$ws.Range("M115").Value = -15050.8568
# Row 116
$ws.Range("H116").Value = 3828.3572
$ws.Range("I116").Value = 3096.1667
$ws.Range("J116").Value = 4377.5
$ws.Range("K116").Value = 3096.1667
$ws.Range("L116").Value = 4377.5
$ws.Range("M116").Value = 345.8332999999998
$ws.Range("N116").Value = -11261.5
# Row 118
$ws.Range("H118").Value = 1660.0
$ws.Range("I118").Value = 1375.0
$ws.Range("K118").Value = 4125.0
$ws.Range("M118").Value = -2468.0
# Row 122
$ws.Range("H122").Value = 731.875
$ws.Range("I122").Value = 309.16666
$ws.Range("K122").Value = 927.4999799999999
$ws.Range("M122").Value = 1522.50002
# Row 125
$ws.Range("H125").Value = 4417.9
$ws.Range("I125").Value = 4078.0
$ws.Range("J125").Value = 4757.8
$ws.Range("K125").Value = 36702.0
$ws.Range("L125").Value = 42820.2
$ws.Range("M125").Value = -34242.0
$ws.Range("N125").Value = -47740.2

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4785.7144
$ws.Range("I32").Value = 4391.9165
$ws.Range("K32").Value = 4391.9165
$ws.Range("M32").Value = -4104.9165
# Row 45
$ws.Range("H45").Value = 2039.8572
$ws.Range("I45").Value = 1505.5454
$ws.Range("K45").Value = 1505.5454
$ws.Range("M45").Value = -1128.5454
# Row 63
$ws.Range("H63").Value = 0.0
$ws.Range("I63").Value = 0.0
$ws.Range("J63").Value = 0.0
$ws.Range("K63").Value = 0.0
$ws.Range("L63").Value = 0.0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 0.0
$ws.Range("I66").Value = 0.0
$ws.Range("J66").Value = 0.0
$ws.Range("K66").Value = 0.0
$ws.Range("L66").Value = 0.0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()
# Row 74
$ws.Range("H74").Value = 11761983.0
$ws.Range("I74").Value = 14281122.0
$ws.Range("K74").Value = 14281122.0
$ws.Range("M74").Value = -14280248.0
# Row 77
$ws.Range("H77").Value = 11761983.0
$ws.Range("I77").Value = 14281122.0
$ws.Range("K77").Value = 71405610.0
$ws.Range("M77").Value = -71401242.0
# Row 88
$ws.Range("H88").Value = 1666.6666
$ws.Range("I88").Value = 1500.5
$ws.Range("J88").Value = 1999.0
$ws.Range("K88").Value = 1500.5
$ws.Range("L88").Value = 1999.0
$ws.Range("M88").Value = -1094.5
$ws.Range("N88").Value = -2811.0
# Row 91
$ws.Range("H91").Value = 1666.6666
$ws.Range("I91").Value = 1500.5
$ws.Range("J91").Value = 1999.0
$ws.Range("K91").Value = 1500.5
$ws.Range("L91").Value = 1999.0
$ws.Range("M91").Value = -96.5
$ws.Range("N91").Value = -4807.0
# Row 122
$ws.Range("H122").Value = 1325.0
$ws.Range("I122").Value = 1404.4546
$ws.Range("J122").Value = 888.0
$ws.Range("K122").Value = 4213.3638
$ws.Range("L122").Value = 2664.0
$ws.Range("M122").Value = -1763.3638
$ws.Range("N122").Value = -7564.0
# Row 132
$ws.Range("H132").Value = 2257.7073
$ws.Range("I132").Value = 1437.9565
$ws.Range("J132").Value = 3305.1667
$ws.Range("K132").Value = 4313.8695
$ws.Range("L132").Value = 9915.500100000001
$ws.Range("M132").Value = -1783.8695
$ws.Range("N132").Value = -14975.5001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1680.3846
$ws.Range("I107").Value = 1571.3334
$ws.Range("K107").Value = 1571.3334
$ws.Range("M107").Value = 348.6666
# Row 115
$ws.Range("H115").Value = 0.0
$ws.Range("J115").Value = 0.0
$ws.Range("L115").Value = 0.0
$ws.Range("N115").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 600.0
$ws.Range("I16").Value = 600.0
$ws.Range("K16").Value = 600.0
$ws.Range("M16").Value = -313.0
# Row 19
$ws.Range("H19").Value = 226.0
$ws.Range("I19").Value = 226.0
$ws.Range("K19").Value = 226.0
$ws.Range("M19").Value = -56.0
# Row 24
$ws.Range("H24").Value = 226.0
$ws.Range("I24").Value = 226.0
$ws.Range("K24").Value = 226.0
$ws.Range("M24").Value = -56.0
# Row 113
$ws.Range("H113").Value = 600.0
$ws.Range("I113").Value = 600.0
$ws.Range("K113").Value = 600.0
$ws.Range("M113").Value = 1570.0
# Row 132
$ws.Range("H132").Value = 2326.182
$ws.Range("I132").Value = 1370.2858
$ws.Range("K132").Value = 4110.857400000001
$ws.Range("M132").Value = -1580.857400000001
# Row 134
$ws.Range("H134").Value = 2680.875
$ws.Range("J134").Value = 1807.0
$ws.Range("L134").Value = 5421.0
$ws.Range("N134").Value = -10491.0

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 31
$ws.Range("H31").Value = 10329.167
$ws.Range("I31").Value = 6395.0
$ws.Range("K31").Value = 6395.0
$ws.Range("M31").Value = -6103.0
# Row 37
$ws.Range("H37").Value = 10329.167
$ws.Range("I37").Value = 6395.0
$ws.Range("K37").Value = 6395.0
$ws.Range("M37").Value = -6118.0
# Row 113
$ws.Range("H113").Value = 2747.25
$ws.Range("I113").Value = 2744.5
$ws.Range("K113").Value = 2744.5
$ws.Range("M113").Value = -574.5
# Row 126
$ws.Range("H126").Value = 0.0
$ws.Range("I126").Value = 0.0
$ws.Range("K126").Value = 0.0
$ws.Range("M126").ClearContents()
# Row 132
$ws.Range("H132").Value = 2943.3333
$ws.Range("I132").Value = 1898.6
$ws.Range("K132").Value = 5695.799999999999
$ws.Range("M132").Value = -3165.799999999999
# Row 139
$ws.Range("H139").Value = 62395.2
$ws.Range("J139").Value = 62395.2
$ws.Range("L139").Value = 62395.2
$ws.Range("N139").Value = -72675.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 5152.3335
$ws.Range("I61").Value = 5152.3335
$ws.Range("K61").Value = 5152.3335
$ws.Range("M61").Value = -4950.3335
# Row 113
$ws.Range("H113").Value = 5152.3335
$ws.Range("I113").Value = 5152.3335
$ws.Range("K113").Value = 5152.3335
$ws.Range("M113").Value = -2982.3335
# Row 132
$ws.Range("H132").Value = 3078.6
$ws.Range("I132").Value = 2131.4
$ws.Range("J132").Value = 4499.4
$ws.Range("K132").Value = 6394.200000000001
$ws.Range("L132").Value = 13498.2
$ws.Range("M132").Value = -3864.200000000001
$ws.Range("N132").Value = -18558.2

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 903.0
$ws.Range("I113").Value = 903.0
$ws.Range("K113").Value = 2709.0
$ws.Range("M113").Value = -539.0
# Row 121
$ws.Range("H121").Value = 49420.0
$ws.Range("J121").Value = 49420.0
$ws.Range("L121").Value = 49420.0
$ws.Range("N121").Value = -52914.0
# Row 123
$ws.Range("H123").Value = 25000.0
$ws.Range("J123").Value = 25000.0
$ws.Range("L123").Value = 25000.0
$ws.Range("N123").Value = -34800.0
# Row 126
$ws.Range("H126").Value = 1319.7646
$ws.Range("I126").Value = 1246.0
$ws.Range("J126").Value = 2500.0
$ws.Range("K126").Value = 3738.0
$ws.Range("L126").Value = 7500.0
$ws.Range("M126").Value = -1268.0
$ws.Range("N126").Value = -12440.0
# Row 132
$ws.Range("H132").Value = 4185.5835
$ws.Range("I132").Value = 4133.5
$ws.Range("J132").Value = 4211.625
$ws.Range("K132").Value = 12400.5
$ws.Range("L132").Value = 12634.875
$ws.Range("M132").Value = -9870.5
$ws.Range("N132").Value = -17694.875

